$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 11.912
$ws.Range("B3").Value = 5.959199999999989
$ws.Range("E5").Value = 12.65789999999999
$ws.Range("B14").Value = 8.934000000000005
$ws.Range("B21").Value = 5.765299999999996
$ws.Range("B23").Value = 5.932599999999995
$ws.Range("B25").Value = 5.943999999999993
